# Updated: st 03. 02. 2021 (New file structure from 07.Feb.2021)
# Refresh the AgTests (H) / AgPosit (I) figures for the existing rows and
# append the new daily row (335) to the Covid daily-stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised AgTests / AgPosit values for previously reported days ---
$ws.Range("H303").Value = 10428
$ws.Range("I303").Value = 669

$ws.Range("H305").Value = 3703
$ws.Range("I305").Value = 304

$ws.Range("H307").Value = 76046
$ws.Range("I307").Value = 6535

$ws.Range("H308").Value = 15727
$ws.Range("I308").Value = 1332

$ws.Range("H309").Value = 75888
$ws.Range("I309").Value = 5299

$ws.Range("H310").Value = 74191
$ws.Range("I310").Value = 3893

$ws.Range("H311").Value = 63018
$ws.Range("I311").Value = 2015

$ws.Range("H312").Value = 26908
$ws.Range("I312").Value = 899

$ws.Range("H313").Value = 61394
$ws.Range("I313").Value = 3269

$ws.Range("H314").Value = 63443
$ws.Range("I314").Value = 3270

$ws.Range("H315").Value = 65792
$ws.Range("I315").Value = 2720

$ws.Range("H316").Value = 49153
$ws.Range("I316").Value = 2225

$ws.Range("H317").Value = 61415
$ws.Range("I317").Value = 2112

$ws.Range("H318").Value = 49056
$ws.Range("I318").Value = 1185

$ws.Range("H319").Value = 41139
$ws.Range("I319").Value = 1627

$ws.Range("H320").Value = 76612
$ws.Range("I320").Value = 3792

$ws.Range("H322").Value = 106206

$ws.Range("H323").Value = 149042

$ws.Range("H324").Value = 230445

$ws.Range("H326").Value = 416620
$ws.Range("I326").Value = 3677

$ws.Range("H327").Value = 236887
$ws.Range("I327").Value = 3537

$ws.Range("H328").Value = 178020
$ws.Range("I328").Value = 2599

$ws.Range("H330").Value = 70609
$ws.Range("I330").Value = 1977

$ws.Range("H331").Value = 146998
$ws.Range("I331").Value = 2536

$ws.Range("H332").Value = 406704
$ws.Range("I332").Value = 3967

$ws.Range("H333").Value = 248133
$ws.Range("I333").Value = 2637

$ws.Range("H334").Value = 197512
$ws.Range("I334").Value = 3381

# --- Append the new daily row (2021-02-07) ---
$ws.Cells.Item(335, 1).Value = 44229
$ws.Cells.Item(335, 2).Value = 254826
$ws.Cells.Item(335, 3).Value = 226471
$ws.Cells.Item(335, 4).Value = 23466
$ws.Cells.Item(335, 5).Value = 12313
$ws.Cells.Item(335, 6).Value = 2732
$ws.Cells.Item(335, 7).Value = 4889
$ws.Cells.Item(335, 8).Value = 111597
$ws.Cells.Item(335, 9).Value = 2592
